# RecruitOn workbook update: filters and maps
# - Adds 4 new applicants (rows 6-9)
# - Replaces the "City, Country" strings with bare city names
# - Adds mailto hyperlinks for the new applicants' D column (Mail)
# - Tweaks Salary / Qualification numbers on a couple of existing rows
# - Moves the active selection to G12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-RowFormat($destRange, $srcRange) {
    $ws.Range($srcRange).Copy()
    $ws.Range($destRange).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# Existing rows 2-5: fix up City + a few numbers
# ---------------------------------------------------------------------------
$ws.Range("F2").Value = "Málaga"

$ws.Range("F3").Value = "Cádiz"
$ws.Range("H3").Value = 16900
$ws.Range("I3").Value = 5.6

$ws.Range("F4").Value = "Madrid"
$ws.Range("H4").Value = 16400
$ws.Range("I4").Value = 4.1

$ws.Range("F5").Value = "Barcelona"
$ws.Range("H5").Value = 19000
$ws.Range("I5").Value = 7.3

# ---------------------------------------------------------------------------
# New rows 6-9: stamp the normal data style (copied from row 2) first, then
# fill in the values so every new cell lands on style "s=2" like its peers.
# ---------------------------------------------------------------------------
Set-RowFormat "A6:I6" "A2:I2"
Set-RowFormat "A7:I7" "A2:I2"
Set-RowFormat "A8:I8" "A2:I2"
Set-RowFormat "A9:I9" "A2:I2"

# Row 6 - Cinco Cincel
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Cinco"
$ws.Range("C6").Value = "Cincel"
$ws.Range("D6").Value = "cinco@cinco.com"
$ws.Range("E6").Value = 600600666
$ws.Range("F6").Value = "Málaga"
$ws.Range("G6").Value = 29
$ws.Range("H6").Value = 12340
$ws.Range("I6").Value = 10

# Row 7 - Seis Sesos
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Seis"
$ws.Range("C7").Value = "Sesos"
$ws.Range("D7").Value = "seis@seis.com"
$ws.Range("E7").Value = 600600665
$ws.Range("F7").Value = "Cádiz"
$ws.Range("G7").Value = 29
$ws.Range("H7").Value = 15700
$ws.Range("I7").Value = 3.3

# Row 8 - Siete Mesino
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Siete"
$ws.Range("C8").Value = "Mesino"
$ws.Range("D8").Value = "siete@siete.com"
$ws.Range("E8").Value = 600600610
$ws.Range("F8").Value = "Madrid"
$ws.Range("G8").Value = 29
$ws.Range("H8").Value = 19230
$ws.Range("I8").Value = 6.34

# Row 9 - Octavo Del Ocho
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Octavo"
$ws.Range("C9").Value = "Del Ocho"
$ws.Range("D9").Value = "oito@orto.com"
$ws.Range("E9").Value = 600600669
$ws.Range("F9").Value = "Barcelona"
$ws.Range("G9").Value = 29
$ws.Range("H9").Value = 24800
$ws.Range("I9").Value = 8.88

# ---------------------------------------------------------------------------
# Hyperlinks on the Mail column for the 4 new rows, styled like D5's link
# (mailto + the yellow-highlight hyperlink look already used on D5).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:cinco@cinco.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:seis@seis.com")
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:siete@siete.com")
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:oito@orto.com")

Set-RowFormat "D6" "D5"
Set-RowFormat "D7" "D5"
Set-RowFormat "D8" "D5"
Set-RowFormat "D9" "D5"

# ---------------------------------------------------------------------------
# Move the selection like the author left it
# ---------------------------------------------------------------------------
$ws.Range("G12").Select()
